$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.291.43'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '1.850.68'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.15'
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6740'
$ws.Range("E6").Value = '  -1.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07461'
$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2974'
$ws.Range("E9").Value = '  -1.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.93'
$ws.Range("E10").Value = '  -0.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07738'
$ws.Range("E11").Value = '  +0.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.036'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").Value = '1.818.10'
$ws.Range("E13").Value = '  -0.88%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6814'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.67'
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.202'
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.187.15'
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008309'
$ws.Range("E18").Value = '  +1.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.59'
$ws.Range("E19").Value = '  +0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.61'
$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9990'
$ws.Range("E21").Value = '  -0.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.255'
$ws.Range("E22").Value = '  -1.95%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9991'
$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.41'
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.722'
$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1411'
$ws.Range("E26").Value = '  -3.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.09'
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.510'
$ws.Range("E28").Value = '  +0.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.203'
$ws.Range("E29").Value = '  -1.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.103'
$ws.Range("E30").Value = '  -1.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.195'
$ws.Range("E31").Value = '  -0.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05345'
$ws.Range("E32").Value = '  +3.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.905'
$ws.Range("E33").Value = '  +3.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7605'
$ws.Range("E34").Value = '  -1.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.147'
$ws.Range("E35").Value = '  +1.18%  '

$ws.Range("E36").Value = '  +0.26%  '

$ws.Range("D37").Value = '1.338.66'
$ws.Range("E37").Value = '  +2.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01813'
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.741'
$ws.Range("E39").Value = '  +1.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9273'
$ws.Range("E40").Value = '  -0.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.971'
$ws.Range("E41").Value = '  +2.93%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.83'
$ws.Range("E43").Value = '  -0.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08026'
$ws.Range("E44").Value = '  +8.34%  '

$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.960.98'
$ws.Range("E45").Value = '  -1.17%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5161'
$ws.Range("E46").Value = '  -0.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.776'
$ws.Range("E47").Value = '  +0.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000124'
$ws.Range("E48").Value = '  +1.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.00'
$ws.Range("E49").Value = '  -3.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.212'
$ws.Range("E50").Value = '  -4.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05941'
$ws.Range("E51").Value = '  +0.38%  '
